$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: Gemini 1.5 Pro results.
$ws.Range("A9").Value = "gemini-1.5-pro"

# Columns with no measured value for this model are still "touched" blank
# cells (matches the source export, which emits an empty cell rather than
# leaving it untouched). Use the leading-apostrophe text marker to force an
# empty text entry, then reset the style so no stray formatting is left
# behind.
foreach ($col in @("B", "E", "F", "G", "H", "I", "O")) {
    $cell = $ws.Range($col + "9")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

$ws.Range("C9").Value = "0.03 ± 0.64"
$ws.Range("D9").Value = "0.33 ± 0.7"
$ws.Range("J9").Value = "0.83 ± 0.08"
$ws.Range("K9").Value = "0.84 ± 0.08"
$ws.Range("L9").Value = "0.84 ± 0.08"
$ws.Range("M9").Value = "0.91 ± 0.1"
$ws.Range("N9").Value = "0.99 ± 0.1"
$ws.Range("P9").Value = "0.58 ± 0.15"
$ws.Range("Q9").Value = "4.83 ± 3.46"
$ws.Range("R9").Value = "0.008 ± 0.00"
$ws.Range("S9").Value = "0.94 ± 0.09"
$ws.Range("T9").Value = "0.9 ± 0.23"
$ws.Range("U9").Value = "2.67 ± 1.28"
$ws.Range("V9").Value = "0.76 ± 0.37"
$ws.Range("W9").Value = "0.94 ± 0.09"
$ws.Range("X9").Value = "1.24 ± 0.29"
